$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '70.867.57'
$ws.Cells.Item(2, 5).Value = '  -0.08%  '

$ws.Range("D3").NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '3.851.51'
$ws.Cells.Item(3, 5).Value = '  +1.77%  '

$ws.Cells.Item(4, 5).Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '698.94'
$ws.Cells.Item(5, 5).Value = '  -0.83%  '

$ws.Range("D6").NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '172.42'
$ws.Cells.Item(6, 5).Value = '  -0.36%  '

$ws.Range("D7").NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '3.848.85'
$ws.Cells.Item(7, 5).Value = '  +1.73%  '

$ws.Cells.Item(8, 5).Value = '  +0.05%  '

$ws.Cells.Item(9, 5).Value = '  -0.07%  '

$ws.Cells.Item(10, 5).Value = '  -0.36%  '

$ws.Range("D11").NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '7.20'
$ws.Cells.Item(11, 5).Value = '  -3.81%  '

$ws.Cells.Item(12, 5).Value = '  -0.64%  '

$ws.Cells.Item(13, 5).Value = '  +0.40%  '

$ws.Range("D14").NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '36.19'
$ws.Cells.Item(14, 5).Value = '  +0.02%  '

$ws.Range("D15").NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '4.501.89'
$ws.Cells.Item(15, 5).Value = '  +1.84%  '

$ws.Range("D16").NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '3.892.30'
$ws.Cells.Item(16, 5).Value = '  +2.84%  '

$ws.Range("D17").NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '70.952.34'
$ws.Cells.Item(17, 5).Value = '  +0.02%  '

$ws.Range("D18").NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '7.17'

$ws.Cells.Item(19, 5).Value = '  -2.86%  '

$ws.Cells.Item(20, 5).Value = '  -0.20%  '

$ws.Range("D21").NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '498.15'
$ws.Cells.Item(21, 5).Value = '  +3.13%  '

$ws.Cells.Item(22, 5).Value = '  -5.17%  '

$ws.Cells.Item(23, 5).Value = '  +0.38%  '

$ws.Cells.Item(24, 5).Value = '  +1.22%  '

$ws.Cells.Item(25, 5).Value = '  +1.57%  '

$ws.Cells.Item(26, 5).Value = '  +0.82%  '

$ws.Range("D27").NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '12.22'
$ws.Cells.Item(27, 5).Value = '  -1.96%  '

$ws.Cells.Item(28, 5).Value = '  -2.78%  '

$ws.Cells.Item(29, 2).Value = 'PancakeSwap'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '3.15'
$ws.Cells.Item(29, 5).Value = '  +0.47%  '

$ws.Cells.Item(30, 2).Value = 'Dai'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D30").NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '0.999'
$ws.Cells.Item(30, 5).Value = '  -0.03%  '

$ws.Range("D31").NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '7.52'
$ws.Cells.Item(31, 5).Value = '  -0.73%  '

$ws.Cells.Item(32, 5).Value = '  -2.81%  '

$ws.Range("D33").NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '29.48'
$ws.Cells.Item(33, 5).Value = '  -0.26%  '

$ws.Cells.Item(34, 5).Value = '  +2.10%  '

$ws.Range("D35").NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '3.808.92'
$ws.Cells.Item(35, 5).Value = '  +1.97%  '

$ws.Range("D36").NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '9.18'
$ws.Cells.Item(36, 5).Value = '  -0.37%  '

$ws.Cells.Item(37, 5).Value = '  -0.04%  '

$ws.Cells.Item(38, 5).Value = '  +0.52%  '

$ws.Cells.Item(39, 5).Value = '  +6.57%  '

$ws.Cells.Item(40, 5).Value = '  +8.21%  '

$ws.Cells.Item(41, 2).Value = 'Filecoin'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '6.01'
$ws.Cells.Item(41, 5).Value = '  +0.14%  '

$ws.Cells.Item(42, 2).Value = 'dogwifhat'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D42").NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '3.37'
$ws.Cells.Item(42, 5).Value = '  -2.69%  '

$ws.Cells.Item(44, 5).Value = '  -0.05%  '

$ws.Range("D45").NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '163.97'
$ws.Cells.Item(45, 5).Value = '  +1.93%  '

$ws.Range("D46").NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '0.000311'
$ws.Cells.Item(46, 5).Value = '  -5.36%  '

$ws.Range("D47").NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '48.98'
$ws.Cells.Item(47, 5).Value = '  -0.32%  '

$ws.Cells.Item(48, 2).Value = 'TheGraph'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D48").NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '0.300'
$ws.Cells.Item(48, 5).Value = '  +0.14%  '

$ws.Cells.Item(49, 2).Value = 'ONDO'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D49").NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '1.38'
$ws.Cells.Item(49, 5).Value = '  -2.94%  '

$ws.Range("D50").NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '43.43'
$ws.Cells.Item(50, 5).Value = '  -6.00%  '

$ws.Cells.Item(51, 5).Value = '  +1.12%  '
